# Fruta / hortaliza, semanal
# Inserts a new weekly price record for Mango at Macroferia Regional de Talca,
# pushing the existing rows 158-179 down to 159-180 (dimension A1:T179 -> A1:T180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 158..179 down by one, making room for the new record at row 158.
$ws.Rows("158:158").Insert()

$ws.Range("A158").Value = 5
$ws.Range("B158").Value = "Macroferia Regional de Talca"
$ws.Range("C158").Value = "Maule"
$ws.Range("D158").Value = 45124
$ws.Range("E158").Value = 7
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100108
$ws.Range("H158").Value = "Tropicales y subtropicales"
$ws.Range("I158").Value = 100108002
$ws.Range("J158").Value = "Mango"
$ws.Range("K158").Value = "Sin especificar"
$ws.Range("L158").Value = "Primera"
$ws.Range("M158").Value = 248
$ws.Range("N158").Value = 8000
$ws.Range("O158").Value = 8000
$ws.Range("P158").Value = 8000
$ws.Range("Q158").Value = "$/bandeja 4 kilos"
$ws.Range("R158").Value = "Brasil"
$ws.Range("S158").Value = 2000
$ws.Range("T158").Value = 4
